$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gains two new columns ("ownTeam", "oppTeam") inserted right after
# "result" (and before the old "batsman" column). Excel's column Insert()
# shifts the existing D:I data (batsman..sr) right to F:K, leaving D:E empty
# for the new columns.
$ws.Columns("D:E").Insert()

# The row contents are also fully reshuffled/extended (7 match rows -> 9),
# so clear everything below the header and rewrite the whole table fresh.
$ws.Range("A2:K8").ClearContents()

# Helper to write a value while always keeping it as TEXT, even for
# cells that look numeric (matches the workbook's existing convention of
# storing every column - including totalRuns/sr/etc - as text, flagged by
# the sheet's numberStoredAsText ignoredError).
function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Header row
$headers = @("venue","date","result","ownTeam","oppTeam","batsman","totalRuns","totalBalls","total4s","total6s","sr")
$cols = @("A","B","C","D","E","F","G","H","I","J","K")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Range($cols[$c] + "1").Value = $headers[$c]
}

# Data rows: venue, date, result, ownTeam, oppTeam, batsman, totalRuns, totalBalls, total4s, total6s, sr
$rows = @(
    @(" Abu Dhabi", " October 28 2020", "Mumbai won by 5 wickets (with 5 balls remaining)", "Royal Challengers Bangalore", "Mumbai Indians", "Washington Sundar ", "10", "6", "1", "0", "166.66"),
    @(" Dubai (DSC)", " September 24 2020", "Kings XI won by 97 runs", "Royal Challengers Bangalore", "Kings XI Punjab", "Washington Sundar ", "30", "27", "2", "1", "111.11"),
    @(" Sharjah", " October 15 2020", "Kings XI won by 8 wickets", "Royal Challengers Bangalore", "Kings XI Punjab", "Washington Sundar ", "13", "14", "1", "0", "92.85"),
    @(" Abu Dhabi", " November 06 2020", "Sunrisers won by 6 wickets (with 2 balls remaining)", "Royal Challengers Bangalore", "Sunrisers Hyderabad", "Washington Sundar ", "5", "6", "0", "0", "83.33"),
    @(" Dubai (DSC)", " October 05 2020", "Capitals won by 59 runs", "Royal Challengers Bangalore", "Delhi Capitals", "Washington Sundar ", "17", "11", "3", "0", "154.54"),
    @(" Sharjah", " October 31 2020", "Sunrisers won by 5 wickets (with 35 balls remaining)", "Royal Challengers Bangalore", "Sunrisers Hyderabad", "Washington Sundar ", "21", "18", "2", "0", "116.66"),
    @(" Abu Dhabi", " November 02 2020", "Capitals won by 6 wickets (with 6 balls remaining)", "Royal Challengers Bangalore", "Delhi Capitals", "Washington Sundar ", "0", "1", "0", "0", "0.00"),
    @(" Dubai (DSC)", " October 25 2020", "Super Kings won by 8 wickets (with 8 balls remaining)", "Royal Challengers Bangalore", "Chennai Super Kings", "Washington Sundar ", "5", "2", "1", "0", "250.00"),
    @(" Dubai (DSC)", " October 10 2020", "RCB won by 37 runs", "Royal Challengers Bangalore", "Chennai Super Kings", "Washington Sundar ", "10", "10", "0", "1", "100.00")
)

# Columns that hold numeric-looking text and therefore need the
# NumberFormat="@" trick so they stay text (index into $cols / row arrays).
$numericCols = @(6, 7, 8, 9, 10)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $addr = $cols[$c] + $rowNum
        if ($numericCols -contains $c) {
            Set-TextCell $addr $row[$c]
        } else {
            $ws.Range($addr).Value = $row[$c]
        }
    }
}
